# Complete rewrite of the internal logic
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header merge: "Autonomous" label moves from F1:K1 to G1:K1,
#     and the "General" label block grows from A1:E1 to A1:F1 ---
$ws.Range("A1:E1").UnMerge()
$ws.Range("F1:K1").UnMerge()

$ws.Range("F1").Value = $null
$ws.Range("G1").Value = "Autonomous"

$ws.Range("A1:F1").Merge()
$ws.Range("G1:K1").Merge()

# Merging resets the cell format on the affected range, so restore the
# original header formatting (border + centered wrap text) by copying it
# from an untouched header cell (L1 keeps its original style).
$ws.Range("L1").Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)

# --- Alliance color vocabulary: "Red" -> "Blue" ---
$ws.Range("E3").Value = "Blue"
$ws.Range("E4").Value = "Blue"
$ws.Range("E5").Value = "Blue"

# --- Row 3 ---
$ws.Range("I3").Value = 2.0
$ws.Range("S3").Value = "Tried but Failed"
$ws.Range("AC3").Value = 0.0
$ws.Range("AD3").Value = 0.0
$ws.Range("AE3").Value = 0.0
$ws.Range("AF3").Value = 0.0

# --- Row 4 ---
$ws.Range("I4").Value = 2.0
$ws.Range("S4").Value = "Tried but Failed"
$ws.Range("AC4").Value = 0.0
$ws.Range("AD4").Value = 0.0
$ws.Range("AE4").Value = 0.0
$ws.Range("AF4").Value = 0.0

# --- Row 5 ---
$ws.Range("B5").Value = "unknown"
$ws.Range("C5").Value = "unknown"
$ws.Range("I5").Value = 0.0
$ws.Range("AE5").Value = 0.0
$ws.Range("AF5").Value = 0.0
